$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RuleSet package cell (C1): was "Auto Doc Gen", now the actual package name.
# Also clear its header styling so it matches the plain/default cell style.
$ws.Range("C1").Value = "com.myspace.eotworkflow"
$ws.Range("C1").Style = "Normal"

# Update the ACTION code cell (F7): parameter name $param -> $1
$ws.Range("F7").Value = "document.setAutogen(`$1)"

# Move the active selection from F6 to C2
[void]$ws.Range("C2").Select()
